$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "65.174.61"
$ws.Cells.Item(2, 5).Value = "  +0.55%  "
$ws.Cells.Item(3, 4).Value = "2.947.76"
$ws.Cells.Item(3, 5).Value = "  -0.93%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.14%  "
$ws.Cells.Item(5, 4).Value = "'567.31"
$ws.Cells.Item(5, 5).Value = "  -2.41%  "
$ws.Cells.Item(6, 4).Value = "'158.22"
$ws.Cells.Item(6, 5).Value = "  +3.27%  "
$ws.Cells.Item(7, 5).Value = "  +0.01%  "
$ws.Cells.Item(8, 4).Value = "'0.520"
$ws.Cells.Item(8, 5).Value = "  +1.32%  "
$ws.Cells.Item(9, 4).Value = "2.945.40"
$ws.Cells.Item(9, 5).Value = "  -0.96%  "
$ws.Cells.Item(10, 5).Value = "  -2.63%  "
$ws.Cells.Item(11, 5).Value = "  +0.76%  "
$ws.Cells.Item(12, 4).Value = "'0.459"
$ws.Cells.Item(12, 5).Value = "  +2.74%  "
$ws.Cells.Item(13, 4).Value = "'0.0000245"
$ws.Cells.Item(13, 5).Value = "  +3.53%  "
$ws.Cells.Item(14, 4).Value = "'34.11"
$ws.Cells.Item(14, 5).Value = "  +0.45%  "
$ws.Cells.Item(15, 5).Value = "  -0.31%  "
$ws.Cells.Item(16, 4).Value = "65.327.87"
$ws.Cells.Item(17, 4).Value = "3.438.35"
$ws.Cells.Item(17, 5).Value = "  -0.87%  "
$ws.Cells.Item(18, 5).Value = "  +1.14%  "
$ws.Cells.Item(19, 4).Value = "2.970.98"
$ws.Cells.Item(19, 5).Value = "  -0.04%  "
$ws.Cells.Item(20, 4).Value = "'447.02"
$ws.Cells.Item(20, 5).Value = "  -0.31%  "
$ws.Cells.Item(21, 4).Value = "'13.90"
$ws.Cells.Item(21, 5).Value = "  +1.66%  "
$ws.Cells.Item(22, 4).Value = "'0.681"
$ws.Cells.Item(22, 5).Value = "  +0.43%  "
$ws.Cells.Item(23, 4).Value = "'7.25"
$ws.Cells.Item(23, 5).Value = "  +0.01%  "
$ws.Cells.Item(24, 4).Value = "'83.05"
$ws.Cells.Item(24, 5).Value = "  +2.64%  "
$ws.Cells.Item(25, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(25, 4).Value = "'12.09"
$ws.Cells.Item(25, 5).Value = "  -1.55%  "
$ws.Cells.Item(26, 2).Value = "Fetch.AI"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(26, 4).Value = "'2.18"
$ws.Cells.Item(26, 5).Value = "  -0.45%  "
$ws.Cells.Item(27, 5).Value = "  -0.01%  "
$ws.Cells.Item(28, 5).Value = "  -6.20%  "
$ws.Cells.Item(29, 4).Value = "'7.93"
$ws.Cells.Item(29, 5).Value = "  +1.48%  "
$ws.Cells.Item(30, 4).Value = "'2.33"
$ws.Cells.Item(30, 5).Value = "  -1.35%  "
$ws.Cells.Item(31, 4).Value = "'2.57"
$ws.Cells.Item(31, 5).Value = "  -0.53%  "
$ws.Cells.Item(32, 4).Value = "0.0₃0984"
$ws.Cells.Item(32, 5).Value = "  -3.12%  "
$ws.Cells.Item(33, 4).Value = "'27.39"
$ws.Cells.Item(33, 5).Value = "  +2.83%  "
$ws.Cells.Item(34, 5).Value = "  +0.39%  "
$ws.Cells.Item(35, 5).Value = "  -0.01%  "
$ws.Cells.Item(36, 4).Value = "'0.974"
$ws.Cells.Item(36, 5).Value = "  -0.73%  "
$ws.Cells.Item(37, 4).Value = "'5.75"
$ws.Cells.Item(37, 5).Value = "  +1.62%  "
$ws.Cells.Item(38, 4).Value = "'49.07"
$ws.Cells.Item(38, 5).Value = "  +0.35%  "
$ws.Cells.Item(40, 4).Value = "'0.299"
$ws.Cells.Item(40, 5).Value = "  +0.85%  "
$ws.Cells.Item(41, 4).Value = "'43.23"
$ws.Cells.Item(41, 5).Value = "  -1.85%  "
$ws.Cells.Item(42, 5).Value = "  -1.10%  "
$ws.Cells.Item(43, 4).Value = "'8.44"
$ws.Cells.Item(43, 5).Value = "  +0.65%  "
$ws.Cells.Item(44, 5).Value = "  -3.66%  "
$ws.Cells.Item(45, 4).Value = "'386.17"
$ws.Cells.Item(45, 5).Value = "  +1.56%  "
$ws.Cells.Item(46, 4).Value = "'0.0353"
$ws.Cells.Item(46, 5).Value = "  +1.59%  "
$ws.Cells.Item(47, 4).Value = "2.746.05"
$ws.Cells.Item(47, 5).Value = "  -0.55%  "
$ws.Cells.Item(48, 4).Value = "'131.78"
$ws.Cells.Item(48, 5).Value = "  -1.68%  "
$ws.Cells.Item(49, 5).Value = "  +0.03%  "
$ws.Cells.Item(50, 2).Value = "ThetaToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(50, 4).Value = "'2.15"
$ws.Cells.Item(50, 5).Value = "  +6.24%  "
$ws.Cells.Item(51, 2).Value = "Stellar"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(51, 4).Value = "'0.107"
$ws.Cells.Item(51, 5).Value = "  +1.53%  "
